$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$quarters = @(
    @(2, "2005Q1"),
    @(3, "2005Q2"),
    @(4, "2005Q3"),
    @(5, "2005Q4"),
    @(6, "2006Q1"),
    @(7, "2006Q2"),
    @(8, "2006Q3"),
    @(9, "2006Q4"),
    @(10, "2007Q1"),
    @(11, "2007Q2"),
    @(12, "2007Q3"),
    @(13, "2007Q4"),
    @(14, "2008Q1"),
    @(15, "2008Q2"),
    @(16, "2008Q3"),
    @(17, "2008Q4"),
    @(18, "2009Q1"),
    @(19, "2009Q2"),
    @(20, "2009Q3"),
    @(21, "2009Q4"),
    @(22, "2010Q1"),
    @(23, "2010Q2"),
    @(24, "2010Q3"),
    @(25, "2010Q4"),
    @(26, "2011Q1"),
    @(27, "2011Q3"),
    @(28, "2011Q4"),
    @(29, "2012Q1"),
    @(30, "2012Q2"),
    @(31, "2012Q3"),
    @(32, "2012Q4"),
    @(33, "2013Q1"),
    @(34, "2013Q2"),
    @(35, "2013Q3"),
    @(36, "2013Q4"),
    @(37, "2014Q1"),
    @(38, "2014Q3"),
    @(39, "2014Q4"),
    @(40, "2015Q1"),
    @(41, "2015Q2"),
    @(42, "2015Q3"),
    @(43, "2015Q4"),
    @(44, "2016Q1"),
    @(45, "2016Q2"),
    @(46, "2016Q3"),
    @(47, "2016Q4"),
    @(48, "2017Q1"),
    @(49, "2017Q2"),
    @(50, "2017Q3"),
    @(51, "2017Q4"),
    @(52, "2018Q1"),
    @(53, "2018Q2"),
    @(54, "2018Q3"),
    @(55, "2018Q4"),
    @(56, "2019Q1"),
    @(57, "2019Q2"),
    @(58, "2019Q3"),
    @(59, "2019Q4"),
    @(60, "2020Q1"),
    @(61, "2020Q2"),
    @(62, "2020Q3"),
    @(63, "2020Q4"),
    @(64, "2021Q1"),
    @(65, "2021Q2"),
    @(66, "2021Q3"),
    @(67, "2021Q4"),
    @(68, "2022Q1"),
    @(69, "2022Q2"),
    @(70, "2022Q3"),
    @(71, "2022Q4"),
    @(72, "2023Q1"),
    @(73, "2023Q2"),
    @(74, "2023Q3"),
    @(75, "2023Q4"),
    @(76, "2024Q1"),
    @(77, "2024Q2"),
    @(78, "2024Q3"),
    @(79, "2024Q4"),
    @(80, "2025Q1"),
    @(81, "2025Q2"),
    @(82, "2025Q3"),
)

$headerStyle = $ws.Cells.Item(1, 1).Style
foreach ($item in $quarters) {
    $r = $item[0]
    $q = $item[1]
    $cell = $ws.Cells.Item($r, 1)
    $cell.Value = $q
    $cell.Style = $headerStyle
}
